$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 460
$ws.Range("F3").Value = 9954
$ws.Range("F4").Value = 215
$ws.Range("F6").Value = 1970
$ws.Range("F7").Value = 6728
$ws.Range("F8").Value = 642
$ws.Range("F9").Value = 0
$ws.Range("F10").Value = 11058
$ws.Range("F11").Value = 12003
$ws.Range("F12").Value = 1276
$ws.Range("F13").Value = 1237
$ws.Range("F14").Value = 5189
$ws.Range("F17").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("F19").Value = 187
$ws.Range("F20").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("F22").Value = 1946
$ws.Range("F23").Value = 955
$ws.Range("F24").Value = 0
$ws.Range("F26").Value = 9
$ws.Range("F28").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("F31").Value = 217
$ws.Range("F33").Value = 0
$ws.Range("F34").Value = 1107
$ws.Range("F35").Value = 0
$ws.Range("F36").Value = 953
$ws.Range("F38").Value = 79
$ws.Range("F39").Value = 3530
$ws.Range("F40").Value = 0
$ws.Range("F42").Value = 556
$ws.Range("F43").Value = 612
$ws.Range("F44").Value = 33
$ws.Range("F45").Value = 0
$ws.Range("F46").Value = 0
$ws.Range("F47").Value = 0
$ws.Range("F48").Value = 0
$ws.Range("F49").Value = 122

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 31
$ws.Range("F6").Value = 10
$ws.Range("F7").Value = 1
$ws.Range("F9").Value = 46
$ws.Range("F10").Value = 1
$ws.Range("F13").Value = 0
$ws.Range("F15").Value = 5
$ws.Range("F18").Value = 0
$ws.Range("F19").Value = 10
$ws.Range("F20").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("F26").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("F29").Value = 0

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6256

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 460
$ws.Range("F3").Value = 9954
$ws.Range("F4").Value = 215
$ws.Range("F6").Value = 13
$ws.Range("F7").Value = 642
$ws.Range("F10").Value = 11058
$ws.Range("F11").Value = 12003
$ws.Range("F12").Value = 33
$ws.Range("F13").Value = 1276
$ws.Range("F14").Value = 0
$ws.Range("F15").Value = 5189
$ws.Range("F16").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("F18").Value = 350
$ws.Range("F19").Value = 46
$ws.Range("F20").Value = 1
$ws.Range("F21").Value = 187
$ws.Range("F22").Value = 72
$ws.Range("F23").Value = 1393
$ws.Range("F24").Value = 304
$ws.Range("F25").Value = 0
$ws.Range("F26").Value = 869
$ws.Range("F27").Value = 468
$ws.Range("F28").Value = 2831
$ws.Range("F29").Value = 0
$ws.Range("F30").Value = 1902
$ws.Range("F31").Value = 100
$ws.Range("F33").Value = 1107
$ws.Range("F34").Value = 0
$ws.Range("F36").Value = 115
$ws.Range("F38").Value = 57
$ws.Range("F40").Value = 46
$ws.Range("F41").Value = 249
$ws.Range("F42").Value = 105
$ws.Range("F43").Value = 556
$ws.Range("F44").Value = 612
$ws.Range("F45").Value = 0
$ws.Range("F48").Value = 0
$ws.Range("F49").Value = 61
